# Generate Report for Handoff
# The handoff markdown file's GUID changed from 690692b4-f56d-4e99-9f57-8d13b074e3fe
# to 6e2d8e7d-b21a-43e0-b7d2-43ad345122f6, a new handoff bundle hash
# (f4fec99773fb3227c418e96dce670ff5cad0323b -> f975771161f660e5f2654a1904ee5fd0b5aac730)
# was generated, and the handoff timestamps were updated for both the zh-cn and
# de-de targets.

$wb = $excel.ActiveWorkbook

$oldGuid = "690692b4-f56d-4e99-9f57-8d13b074e3fe"
$newGuid = "6e2d8e7d-b21a-43e0-b7d2-43ad345122f6"
$oldHash = "f4fec99773fb3227c418e96dce670ff5cad0323b"
$newHash = "f975771161f660e5f2654a1904ee5fd0b5aac730"

$mdName = "$newGuid.md"
$zhCnXlfName = "$newGuid.$newHash.zh-cn.xlf"
$deDeXlfName = "$newGuid.$newHash.de-de.xlf"

$zhCnDatetime = "2016-03-09 16:51:00"
$deDeDatetime = "2016-03-09 16:51:06"

# ---------------------------------------------------------------------------
# Overview sheet: just the file-name cell/hyperlink in A2
# ---------------------------------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")

$wsOverview.Cells.Item(2, 1).Value = $mdName

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(2, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/fb54fccd1599754e243217ff0df14bbe3b2ea1b8/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/fb54fccd1599754e243217ff0df14bbe3b2ea1b8/.localization-config",
    "",
    "",
    ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: file-name in A2, handoff file name + datetime in C2/D2
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Sheets.Item("zh-cn")

$wsZhCn.Cells.Item(2, 1).Value = $mdName
$wsZhCn.Cells.Item(2, 3).Value = $zhCnXlfName
$wsZhCn.Cells.Item(2, 4).Value = $zhCnDatetime

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(2, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/fb54fccd1599754e243217ff0df14bbe3b2ea1b8/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(2, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/44f5f8c03b2c66eb8347d662f1a20e8598643dff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhCnXlfName",
    "",
    "",
    $zhCnXlfName) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/fb54fccd1599754e243217ff0df14bbe3b2ea1b8/.localization-config",
    "",
    "",
    ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: file-name in A2, handoff file name + datetime in C2/D2
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Sheets.Item("de-de")

$wsDeDe.Cells.Item(2, 1).Value = $mdName
$wsDeDe.Cells.Item(2, 3).Value = $deDeXlfName
$wsDeDe.Cells.Item(2, 4).Value = $deDeDatetime

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(2, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/fb54fccd1599754e243217ff0df14bbe3b2ea1b8/e2e/$mdName",
    "",
    "",
    $mdName) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(2, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a9146784acd6d528f5cf16d35b7893ae568b91a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deDeXlfName",
    "",
    "",
    $deDeXlfName) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/fb54fccd1599754e243217ff0df14bbe3b2ea1b8/.localization-config",
    "",
    "",
    ".localization-config") | Out-Null
